$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '28.606.13'
Set-TextValue $ws.Range("E2") '  +0.97%  '
Set-TextValue $ws.Range("D3") '1.564.24'
Set-TextValue $ws.Range("E3") '  +0.19%  '
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  -0.13%  '
Set-TextValue $ws.Range("D5") '210.29'
Set-TextValue $ws.Range("E5") '  -0.25%  '
Set-TextValue $ws.Range("D6") '0.505'
Set-TextValue $ws.Range("E6") '  +3.09%  '
Set-TextValue $ws.Range("E7") '  -0.16%  '
Set-TextValue $ws.Range("D8") '24.88'
Set-TextValue $ws.Range("E8") '  +5.62%  '
Set-TextValue $ws.Range("E9") '  +0.79%  '
Set-TextValue $ws.Range("E10") '  -0.04%  '
Set-TextValue $ws.Range("D11") '0.0896'
Set-TextValue $ws.Range("E11") '  +0.42%  '
Set-TextValue $ws.Range("D12") '1.787.61'
Set-TextValue $ws.Range("E12") '  +0.21%  '
Set-TextValue $ws.Range("D13") '1.567.93'
Set-TextValue $ws.Range("E13") '  +0.51%  '
Set-TextValue $ws.Range("D14") '28.640.73'
Set-TextValue $ws.Range("E14") '  +1.14%  '
Set-TextValue $ws.Range("E15") '  +0.75%  '
Set-TextValue $ws.Range("E16") '  -0.59%  '
Set-TextValue $ws.Range("D17") '61.36'
Set-TextValue $ws.Range("E17") '  +0.65%  '
Set-TextValue $ws.Range("D18") '227.40'
Set-TextValue $ws.Range("E18") '  -0.16%  '
Set-TextValue $ws.Range("E19") '  -0.48%  '
Set-TextValue $ws.Range("D20") '0.0₃0680'
Set-TextValue $ws.Range("E20") '  +0.34%  '
Set-TextValue $ws.Range("D21") '1.00'
Set-TextValue $ws.Range("E21") '  -0.13%  '
Set-TextValue $ws.Range("D22") '3.93'
Set-TextValue $ws.Range("E22") '  -0.05%  '
Set-TextValue $ws.Range("D23") '9.04'
Set-TextValue $ws.Range("E23") '  +1.41%  '
Set-TextValue $ws.Range("E24") '  +1.00%  '
Set-TextValue $ws.Range("D25") '151.84'
Set-TextValue $ws.Range("E25") '  +1.01%  '
Set-TextValue $ws.Range("D26") '0.105'
Set-TextValue $ws.Range("E26") '  +1.89%  '
Set-TextValue $ws.Range("D27") '14.78'
Set-TextValue $ws.Range("E28") '  -0.11%  '
Set-TextValue $ws.Range("D29") '6.24'
Set-TextValue $ws.Range("E29") '  -1.46%  '
Set-TextValue $ws.Range("D30") '0.0459'
Set-TextValue $ws.Range("E30") '  -3.64%  '
Set-TextValue $ws.Range("E31") '  -0.50%  '
Set-TextValue $ws.Range("E32") '  +0.06%  '
Set-TextValue $ws.Range("D33") '1.404.56'
Set-TextValue $ws.Range("E33") '  +1.52%  '
Set-TextValue $ws.Range("E34") '  -2.72%  '
Set-TextValue $ws.Range("D35") '1.03'
Set-TextValue $ws.Range("E35") '  -3.05%  '
Set-TextValue $ws.Range("D36") '1.47'
Set-TextValue $ws.Range("E36") '  -1.65%  '
Set-TextValue $ws.Range("E37") '  +1.43%  '
Set-TextValue $ws.Range("E38") '  -1.86%  '
Set-TextValue $ws.Range("D39") '0.0162'
Set-TextValue $ws.Range("E39") '  -0.36%  '
Set-TextValue $ws.Range("E40") '  +0.12%  '
Set-TextValue $ws.Range("D41") '0.516'
Set-TextValue $ws.Range("E41") '  -0.52%  '
Set-TextValue $ws.Range("D42") '1.00'
Set-TextValue $ws.Range("E42") '  -0.12%  '
Set-TextValue $ws.Range("D43") '0.767'
Set-TextValue $ws.Range("E43") '  -1.89%  '
Set-TextValue $ws.Range("D44") '0.0461'
Set-TextValue $ws.Range("E44") '  -2.07%  '
Set-TextValue $ws.Range("D45") '63.88'
Set-TextValue $ws.Range("E45") '  +2.73%  '
Set-TextValue $ws.Range("E46") '  -2.14%  '
Set-TextValue $ws.Range("D47") '1.699.25'
Set-TextValue $ws.Range("D48") '0.865'
Set-TextValue $ws.Range("E48") '  -5.67%  '
Set-TextValue $ws.Range("D49") '84.90'
Set-TextValue $ws.Range("E49") '  -0.61%  '
Set-TextValue $ws.Range("E50") '  +4.72%  '
Set-TextValue $ws.Range("D51") '0.0510'
Set-TextValue $ws.Range("E51") '  -0.49%  '
